# Integrated Recomendation Engine with REST API
#
# Updates the major-recommendation matrix on Sheet1:
#  - Row 11 (MSE 200 ...): Physics is now a secondary major tag.
#  - Row 12 (MAE 416 ...): Physics is now a secondary major tag, and the
#    Physics column for this row is re-targeted to a Physical
#    Chemistry / Physics+Chemistry pairing.
#  - Minor pluralization fix: "Humanities, Social Sciences" -> the
#    canonical "Humanities, Social Science" tag used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J12").Value = "Physics, Chemistry"
$ws.Range("I12").Value = "Physical Chemistry"

$ws.Range("B11").Value = "Engineering, Physics"
$ws.Range("F11").Value = "Humanities, Social Science"
$ws.Range("B12").Value = "Engineering, Physics"

# Restore the active selection/scroll position to match the author's
# final view of the sheet (selection moved from P10 to B16, and the
# window scrolled back to show column A).
$ws.Activate() | Out-Null
$ws.Range("B16").Select() | Out-Null
